# Updates cryptos list data cells (Coin/Link/Price/Volume(1h)) on Sheet1
# to match the refreshed GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.208.04"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "3.589.81"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.61"
$ws.Range("E5").Value = "  -2.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.46"
$ws.Range("E6").Value = "  -1.44%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  -2.50%  "

$ws.Range("D8").Value = "3.587.47"
$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("E10").Value = "  -0.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.662"
$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.93"
$ws.Range("E12").Value = "  -3.19%  "

$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.64"
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("D15").Value = "4.168.22"
$ws.Range("E15").Value = "  -1.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.95"
$ws.Range("E16").Value = "  +3.10%  "

$ws.Range("D17").Value = "3.588.80"
$ws.Range("E17").Value = "  -1.33%  "

$ws.Range("D18").Value = "70.101.66"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.64"
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("E21").Value = "  -0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "476.03"
$ws.Range("E22").Value = "  -3.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.13"
$ws.Range("E23").Value = "  +14.87%  "

$ws.Range("E24").Value = "  -8.09%  "

$ws.Range("E25").Value = "  -1.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.16"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("E27").Value = "  -1.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.05"
$ws.Range("E28").Value = "  -1.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.25"
$ws.Range("E30").Value = "  -0.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.72"
$ws.Range("E31").Value = "  +1.58%  "

$ws.Range("E32").Value = "  +3.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.17"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.32"
$ws.Range("E34").Value = "  +1.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "583.63"
$ws.Range("E35").Value = "  -5.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.11"
$ws.Range("E36").Value = "  +2.96%  "

$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").Value = "0.0₃0797"
$ws.Range("E38").Value = "  -4.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.397"
$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.27"
$ws.Range("E40").Value = "  +18.75%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.140"
$ws.Range("E41").Value = "  -5.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.49"
$ws.Range("E42").Value = "  -4.99%  "

$ws.Range("D43").Value = "3.244.05"
$ws.Range("E43").Value = "  -3.42%  "

$ws.Range("E44").Value = "  +7.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.10"
$ws.Range("E45").Value = "  +1.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0444"
$ws.Range("E46").Value = "  -0.22%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.34"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.43"
$ws.Range("E48").Value = "  +3.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.138"
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.16"
$ws.Range("E51").Value = "  -3.86%  "
